$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94; this shifts the existing rows 94-140
# down to 95-141, carrying their values/formatting with them.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new record's data.
$row = 94
$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).Value = 45205
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(95, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 100112037
$ws.Cells.Item($row, 7).Value = "Cebollín"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 70
$ws.Cells.Item($row, 11).Value = 4000
$ws.Cells.Item($row, 12).Value = 4500
$ws.Cells.Item($row, 13).Value = 4214
$ws.Cells.Item($row, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 117
$ws.Cells.Item($row, 17).Value = 36
$ws.Cells.Item($row, 18).Value = "Hortaliza"
